# Generate Report for Handback
# - Update status text from "Ready for handoff" to "Handed back: in sync with en-US"
#   on every sheet that shows it (Overview, zh-cn, de-de).
# - Refresh the "Latest Handback DateTime" timestamps for zh-cn / de-de.
# - Clear the stale "handback file is not latest" Error Detail messages now
#   that the handback is in sync.
# - Widen the Status / zh-cn / de-de columns and narrow the now-empty Error
#   Detail column to match the regenerated report layout.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

$overview.Columns.Item(5).ColumnWidth = 29.14437166849777
$overview.Columns.Item(6).ColumnWidth = 29.14437166849777

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$zhcn.Range("K2").Value = "2016-09-05 10:43:37"
$zhcn.Range("K3").Value = "2016-09-05 10:43:37"

$zhcn.Range("P2").Value = ""
$zhcn.Range("P3").Value = ""

$zhcn.Columns.Item(3).ColumnWidth = 29.14437166849777
$zhcn.Columns.Item(16).ColumnWidth = 12.913719813028965

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Range("K2").Value = "2016-09-05 10:43:57"
$dede.Range("K3").Value = "2016-09-05 10:43:57"

$dede.Range("P2").Value = ""
$dede.Range("P3").Value = ""

$dede.Columns.Item(3).ColumnWidth = 29.14437166849777
$dede.Columns.Item(16).ColumnWidth = 12.913719813028965
